# Generate Report for Handoff
#
# The localization-status report is re-sorted: the row for
# "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" (previously "Handed back: in
# sync with en-US") moves from the 2nd data row to the last data row, its
# Status becomes "Ready for handoff", and its handoff timestamp advances.
# The other two rows shift up to fill the gap. Hyperlink display text is
# kept in sync with the cell it decorates.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value2 = "ffff27369140-abc4-4010-88af-940ba64a7fee.md"
$ws1.Range("A3").Value2 = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md"
$ws1.Range("A4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"

foreach ($hl in $ws1.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff27369140-abc4-4010-88af-940ba64a7fee.md" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" }
}

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value2 = "ffff27369140-abc4-4010-88af-940ba64a7fee.md"
$ws2.Range("C2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.zh-cn.xlf"
$ws2.Range("D2").Value2 = "2016-03-04 06:24:44"
$ws2.Range("E2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.md"
$ws2.Range("F2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.zh-cn.xlf"
$ws2.Range("G2").Value2 = "2016-03-04 06:25:51"

$ws2.Range("A3").Value2 = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md"

$ws2.Range("A4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md"
$ws2.Range("B4").Value2 = "Ready for handoff"
$ws2.Range("C4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.zh-cn.xlf"
$ws2.Range("D4").Value2 = "2016-03-04 06:30:24"
$ws2.Range("E4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md"
$ws2.Range("F4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.zh-cn.xlf"
$ws2.Range("G4").Value2 = "2016-03-04 06:29:06"

foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff27369140-abc4-4010-88af-940ba64a7fee.md" }
    elseif ($addr -eq '$C$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.zh-cn.xlf" }
    elseif ($addr -eq '$E$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.md" }
    elseif ($addr -eq '$F$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.zh-cn.xlf" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.zh-cn.xlf" }
    elseif ($addr -eq '$E$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" }
    elseif ($addr -eq '$F$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.zh-cn.xlf" }
}

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value2 = "ffff27369140-abc4-4010-88af-940ba64a7fee.md"
$ws3.Range("C2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.de-de.xlf"
$ws3.Range("D2").Value2 = "2016-03-04 06:24:59"
$ws3.Range("E2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.md"
$ws3.Range("F2").Value2 = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.de-de.xlf"
$ws3.Range("G2").Value2 = "2016-03-04 06:26:21"

$ws3.Range("A3").Value2 = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md"

$ws3.Range("A4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md"
$ws3.Range("B4").Value2 = "Ready for handoff"
$ws3.Range("C4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.de-de.xlf"
$ws3.Range("D4").Value2 = "2016-03-04 06:30:40"
$ws3.Range("E4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md"
$ws3.Range("F4").Value2 = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.de-de.xlf"
$ws3.Range("G4").Value2 = "2016-03-04 06:29:34"

foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff27369140-abc4-4010-88af-940ba64a7fee.md" }
    elseif ($addr -eq '$C$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.de-de.xlf" }
    elseif ($addr -eq '$E$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.md" }
    elseif ($addr -eq '$F$2') { $hl.TextToDisplay = "465496d9-901d-4e13-a66e-e96712b17117.5686d97c36b38c165bf6378a8fbc03f2a531c9f4.de-de.xlf" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffffa1d25d47-cea4-45a7-963a-ae0928b18d1e.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.de-de.xlf" }
    elseif ($addr -eq '$E$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.md" }
    elseif ($addr -eq '$F$4') { $hl.TextToDisplay = "03741d24-08a6-4b4a-82ee-ef40b1d66af5.a10bd74dbf3ab657d3edde4b53a6448ddaeeb5eb.de-de.xlf" }
}

Write-Output "Done updating handoff report."
